{"js": "// Update the worksheet date (first paragraph in the document body).\nconst body = context.document.body;\n\nconst dateResults = body.search(\"2024-03-16 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-03-17 Sunday\", Word.InsertLocation.replace);\n}\n\n// Update the 20x5 grid of addition/subtraction problems in the single table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New problem text, row-major (20 rows x 5 columns), matching the table's\n// current layout exactly so each cell keeps its existing formatting.\ntable.values = [\n  [\"67-56=\", \"15-8=\", \"95-69=\", \"76+21=\", \"80+5=\"],\n  [\"1+53=\", \"41+22=\", \"40-26=\", \"58-15=\", \"4+0=\"],\n  [\"20+33=\", \"44+31=\", \"3+13=\", \"16+7=\", \"98-66=\"],\n  [\"70-37=\", \"50-8=\", \"70-61=\", \"29-12=\", \"24+39=\"],\n  [\"61-23=\", \"2+70=\", \"69+20=\", \"7+47=\", \"45-27=\"],\n  [\"64+29=\", \"41+12=\", \"34+47=\", \"83-54=\", \"21+30=\"],\n  [\"56-9=\", \"50-43=\", \"90+4=\", \"51+36=\", \"44+52=\"],\n  [\"72-10=\", \"36+34=\", \"96-78=\", \"44+0=\", \"98-25=\"],\n  [\"77-18=\", \"8+31=\", \"4+2=\", \"62-56=\", \"26+9=\"],\n  [\"43-35=\", \"46+12=\", \"75-70=\", \"84-76=\", \"18+6=\"],\n  [\"99-23=\", \"76-14=\", \"70-65=\", \"98-44=\", \"84-64=\"],\n  [\"35-3=\", \"58-52=\", \"8+48=\", \"13+41=\", \"28-6=\"],\n  [\"2+65=\", \"45+6=\", \"50-2=\", \"97-4=\", \"75+15=\"],\n  [\"55+9=\", \"16+44=\", \"82-74=\", \"17+14=\", \"21+60=\"],\n  [\"30+18=\", \"18-14=\", \"62-58=\", \"34+16=\", \"66-18=\"],\n  [\"11+55=\", \"43-23=\", \"48+28=\", \"52+8=\", \"50-22=\"],\n  [\"64-46=\", \"10+4=\", \"32-26=\", \"26-9=\", \"29+64=\"],\n  [\"25+70=\", \"31+30=\", \"72-17=\", \"9+52=\", \"31-29=\"],\n  [\"78+5=\", \"90-68=\", \"79-21=\", \"58-0=\", \"62+32=\"],\n  [\"38+23=\", \"99-72=\", \"44+20=\", \"61-25=\", \"23+65=\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the worksheet date (first paragraph in the document).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"2024-03-16 Saturday\",   # FindText\n    $false,                  # MatchCase\n    $false,                  # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                  # Format\n    \"2024-03-17 Sunday\",     # ReplaceWith\n    2                        # Replace (wdReplaceAll)\n)\n\n# Update the 20x5 grid of addition/subtraction problems in the single table.\n$newValues = @(\n    @(\"67-56=\", \"15-8=\", \"95-69=\", \"76+21=\", \"80+5=\"),\n    @(\"1+53=\", \"41+22=\", \"40-26=\", \"58-15=\", \"4+0=\"),\n    @(\"20+33=\", \"44+31=\", \"3+13=\", \"16+7=\", \"98-66=\"),\n    @(\"70-37=\", \"50-8=\", \"70-61=\", \"29-12=\", \"24+39=\"),\n    @(\"61-23=\", \"2+70=\", \"69+20=\", \"7+47=\", \"45-27=\"),\n    @(\"64+29=\", \"41+12=\", \"34+47=\", \"83-54=\", \"21+30=\"),\n    @(\"56-9=\", \"50-43=\", \"90+4=\", \"51+36=\", \"44+52=\"),\n    @(\"72-10=\", \"36+34=\", \"96-78=\", \"44+0=\", \"98-25=\"),\n    @(\"77-18=\", \"8+31=\", \"4+2=\", \"62-56=\", \"26+9=\"),\n    @(\"43-35=\", \"46+12=\", \"75-70=\", \"84-76=\", \"18+6=\"),\n    @(\"99-23=\", \"76-14=\", \"70-65=\", \"98-44=\", \"84-64=\"),\n    @(\"35-3=\", \"58-52=\", \"8+48=\", \"13+41=\", \"28-6=\"),\n    @(\"2+65=\", \"45+6=\", \"50-2=\", \"97-4=\", \"75+15=\"),\n    @(\"55+9=\", \"16+44=\", \"82-74=\", \"17+14=\", \"21+60=\"),\n    @(\"30+18=\", \"18-14=\", \"62-58=\", \"34+16=\", \"66-18=\"),\n    @(\"11+55=\", \"43-23=\", \"48+28=\", \"52+8=\", \"50-22=\"),\n    @(\"64-46=\", \"10+4=\", \"32-26=\", \"26-9=\", \"29+64=\"),\n    @(\"25+70=\", \"31+30=\", \"72-17=\", \"9+52=\", \"31-29=\"),\n    @(\"78+5=\", \"90-68=\", \"79-21=\", \"58-0=\", \"62+32=\"),\n    @(\"38+23=\", \"99-72=\", \"44+20=\", \"61-25=\", \"23+65=\")\n)\n\n$table = $d.Tables(1)\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
